# Apply updated "F" column values to the "展览" and "全部类型" worksheets,
# matching the commit's regenerated data snapshot.

$wb = $excel.ActiveWorkbook

$wsExhibition = $wb.Worksheets.Item("展览")
$wsExhibition.Range("F2").Value = 1088
$wsExhibition.Range("F3").Value = 398
$wsExhibition.Range("F4").Value = 1499
$wsExhibition.Range("F5").Value = 8748
$wsExhibition.Range("F6").Value = 94
$wsExhibition.Range("F10").Value = 152
$wsExhibition.Range("F11").Value = 17
$wsExhibition.Range("F12").Value = 3596
$wsExhibition.Range("F15").Value = 79
$wsExhibition.Range("F16").Value = 1318
$wsExhibition.Range("F20").Value = 204
$wsExhibition.Range("F21").Value = 2365

$wsAllTypes = $wb.Worksheets.Item("全部类型")
$wsAllTypes.Range("F2").Value = 1088
$wsAllTypes.Range("F3").Value = 398
$wsAllTypes.Range("F4").Value = 1499
$wsAllTypes.Range("F5").Value = 8748
$wsAllTypes.Range("F6").Value = 94
$wsAllTypes.Range("F10").Value = 152
$wsAllTypes.Range("F11").Value = 17
$wsAllTypes.Range("F12").Value = 3596
$wsAllTypes.Range("F15").Value = 79
$wsAllTypes.Range("F16").Value = 1319
$wsAllTypes.Range("F20").Value = 204
$wsAllTypes.Range("F21").Value = 2365
